# Applies the "456a3b4" data refresh to 杭州-漫展信息.xlsx
#  - Sheet "展览"   (Exhibitions):  bump "想去人数" (F column) counters
#  - Sheet "本地生活" (Local life):  append a newly-scraped event as row 2
#  - Sheet "全部类型" (All types):   bump the same "想去人数" counters
# (Sheet "演出" / Performances is untouched.)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------------
# 1. 展览 ("Exhibitions") sheet - updated "想去人数" (column F) values
# ---------------------------------------------------------------------
$ws1.Cells.Item(4,6).Value  = 320
$ws1.Cells.Item(6,6).Value  = 1161
$ws1.Cells.Item(8,6).Value  = 58
$ws1.Cells.Item(11,6).Value = 1225
$ws1.Cells.Item(14,6).Value = 910
$ws1.Cells.Item(18,6).Value = 79
$ws1.Cells.Item(20,6).Value = 813
$ws1.Cells.Item(21,6).Value = 1755
$ws1.Cells.Item(22,6).Value = 3153
$ws1.Cells.Item(23,6).Value = 922
$ws1.Cells.Item(25,6).Value = 2307
$ws1.Cells.Item(26,6).Value = 672
$ws1.Cells.Item(27,6).Value = 10
$ws1.Cells.Item(28,6).Value = 3173
$ws1.Cells.Item(29,6).Value = 652
$ws1.Cells.Item(30,6).Value = 807
$ws1.Cells.Item(31,6).Value = 20
$ws1.Cells.Item(32,6).Value = 94
$ws1.Cells.Item(33,6).Value = 745
$ws1.Cells.Item(34,6).Value = 151
$ws1.Cells.Item(35,6).Value = 142
$ws1.Cells.Item(36,6).Value = 63
$ws1.Cells.Item(38,6).Value = 1130
$ws1.Cells.Item(39,6).Value = 1815
$ws1.Cells.Item(40,6).Value = 415
$ws1.Cells.Item(43,6).Value = 210
$ws1.Cells.Item(44,6).Value = 139
$ws1.Cells.Item(45,6).Value = 189

# ---------------------------------------------------------------------
# 2. 本地生活 ("Local life") sheet - insert the newly scraped event as
#    row 2 (dimension grows from A1:I1 to A1:I2)
# ---------------------------------------------------------------------

# A2 should inherit the same bold/centered/bordered style as A1 (s="1"),
# so copy formatting only, then fill in the actual values/text below.
$ws3.Range("A1").Copy() | Out-Null
$ws3.Range("A2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

$ws3.Cells.Item(2,1).Value = 1

# Columns B and E contain free-form date-like text ("2024-04-19",
# "2024.04.19 00:00-05.10 23:59") that must stay literal text rather
# than being auto-converted into a date serial number.
$ws3.Range("B2").NumberFormat = "@"
$ws3.Range("B2").Value = "2024-04-19"
$ws3.Range("B2").Style = "Normal"

$ws3.Cells.Item(2,3).Value = "杭州·偶像梦幻祭2「绽放的纯白浪漫」线下特别快闪"
$ws3.Cells.Item(2,4).Value = "延安路292号（地铁1号线龙翔桥站D出口） 工联CC"

$ws3.Range("E2").NumberFormat = "@"
$ws3.Range("E2").Value = "2024.04.19 00:00-05.10 23:59"
$ws3.Range("E2").Style = "Normal"

$ws3.Cells.Item(2,6).Value = 75
$ws3.Cells.Item(2,7).Value = 10
$ws3.Cells.Item(2,8).Value = "https://show.bilibili.com/platform/detail.html?id=84042"
$ws3.Cells.Item(2,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/DkJm2r5z1712745086412.jpeg"

# ---------------------------------------------------------------------
# 3. 全部类型 ("All types") sheet - same "想去人数" (column F) bumps as
#    in 展览, but the row numbers differ because this sheet aggregates
#    rows from every other sheet.
# ---------------------------------------------------------------------
$ws4.Cells.Item(4,6).Value  = 320
$ws4.Cells.Item(5,6).Value  = 1161
$ws4.Cells.Item(8,6).Value  = 1225
$ws4.Cells.Item(10,6).Value = 910
$ws4.Cells.Item(16,6).Value = 79
$ws4.Cells.Item(17,6).Value = 813
$ws4.Cells.Item(18,6).Value = 1755
$ws4.Cells.Item(19,6).Value = 3153
$ws4.Cells.Item(20,6).Value = 922
$ws4.Cells.Item(23,6).Value = 2307
$ws4.Cells.Item(24,6).Value = 10
$ws4.Cells.Item(25,6).Value = 3173
$ws4.Cells.Item(26,6).Value = 652
$ws4.Cells.Item(27,6).Value = 807
$ws4.Cells.Item(29,6).Value = 20
$ws4.Cells.Item(33,6).Value = 94
$ws4.Cells.Item(35,6).Value = 745
$ws4.Cells.Item(36,6).Value = 151
$ws4.Cells.Item(37,6).Value = 142
$ws4.Cells.Item(38,6).Value = 63
$ws4.Cells.Item(41,6).Value = 1130
$ws4.Cells.Item(42,6).Value = 1815
$ws4.Cells.Item(44,6).Value = 415
$ws4.Cells.Item(46,6).Value = 210
$ws4.Cells.Item(47,6).Value = 139
$ws4.Cells.Item(48,6).Value = 189
